$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("E1").Value = "Number of levels in game"
$ws.Range("F1").Value = "Category"
$ws.Range("G1").Value = "Level From"
$ws.Range("H1").Value = "Level To"

# --- New data cells added alongside the existing rows ---
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 5

$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = 13

$ws.Range("E15").Value = 4

$ws.Range("E16").Value = 1

# --- Column widths for the newly used columns ---
$ws.Columns.Item(4).ColumnWidth = 28.833333333333336
$ws.Columns.Item(5).ColumnWidth = 24.833333333333336
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666
$ws.Columns.Item(7).ColumnWidth = 12.5

# --- Selection / view state ---
$ws.Range("L13").Select() | Out-Null
